$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.029.91'
$ws.Range('E2').Value = '  -0.66%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.421.70'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '554.20'
$ws.Range('E5').Value = '  -0.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.50'
$ws.Range('E6').Value = '  -0.81%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.593'
$ws.Range('E8').Value = '  +3.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.106'
$ws.Range('E9').Value = '  -1.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.67'
$ws.Range('E10').Value = '  -2.46%  '
$ws.Range('E11').Value = '  -0.84%  '
$ws.Range('E12').Value = '  -1.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '25.29'
$ws.Range('E13').Value = '  +1.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.855.65'
$ws.Range('E14').Value = '  -0.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.008.75'
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('E16').Value = '  -1.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.425.98'
$ws.Range('E17').Value = '  +0.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.30'
$ws.Range('E18').Value = '  -1.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.46'
$ws.Range('E19').Value = '  +1.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '328.56'
$ws.Range('E20').Value = '  -2.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.68'
$ws.Range('E21').Value = '  -2.72%  '
$ws.Range('E22').Value = '  +0.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.02'
$ws.Range('E23').Value = '  +1.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.179'
$ws.Range('E24').Value = '  +4.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.61'
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.41'
$ws.Range('E27').Value = '  +3.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0778'
$ws.Range('E28').Value = '  -0.83%  '
$ws.Range('E29').Value = '  -2.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '169.47'
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.08'
$ws.Range('E31').Value = '  -3.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.08'
$ws.Range('E32').Value = '  +3.68%  '
$ws.Range('B33').Value = 'PolygonEcosystemToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.404'
$ws.Range('E33').Value = '  -3.76%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.63'
$ws.Range('E34').Value = '  -0.99%  '
$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.32'
$ws.Range('E36').Value = '  +0.47%  '
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.20'
$ws.Range('E38').Value = '  -1.11%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '327.26'
$ws.Range('E39').Value = '  +3.01%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.60'
$ws.Range('E40').Value = '  -1.68%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '141.41'
$ws.Range('E41').Value = '  -0.97%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.66'
$ws.Range('E42').Value = '  -1.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0971'
$ws.Range('E43').Value = '  +1.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.63'
$ws.Range('E44').Value = '  -1.22%  '
$ws.Range('E45').Value = '  -1.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.575'
$ws.Range('E46').Value = '  +0.41%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0224'
$ws.Range('E47').Value = '  -1.26%  '
$ws.Range('B48').Value = 'Polygon'
$ws.Range('C48').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.387'
$ws.Range('E48').Value = '  -5.48%  '
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('E50').Value = '  -2.85%  '
$ws.Range('E51').Value = '  -1.08%  '
